$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for most rows
$ws.Range("D2").Value = "69.269.18"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "3.672.73"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "680.28"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").Value = "158.96"
$ws.Range("E6").Value = "  -7.23%  "
$ws.Range("D7").Value = "3.669.61"
$ws.Range("E7").Value = "  -3.59%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  -6.34%  "
$ws.Range("E10").Value = "  -10.04%  "
$ws.Range("D11").Value = "7.08"
$ws.Range("E11").Value = "  -5.14%  "
$ws.Range("D12").Value = "0.433"
$ws.Range("E12").Value = "  -10.21%  "
$ws.Range("E13").Value = "  -7.85%  "
$ws.Range("D14").Value = "4.291.41"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").Value = "32.25"
$ws.Range("E15").Value = "  -11.29%  "
$ws.Range("D16").Value = "3.679.29"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "69.242.36"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "15.77"
$ws.Range("E19").Value = "  -9.90%  "
$ws.Range("E20").Value = "  -11.40%  "
$ws.Range("D21").Value = "471.81"
$ws.Range("E21").Value = "  -8.59%  "
$ws.Range("D22").Value = "9.83"
$ws.Range("E22").Value = "  -5.89%  "
$ws.Range("E23").Value = "  -9.67%  "
$ws.Range("D24").Value = "79.15"
$ws.Range("E24").Value = "  -5.83%  "
$ws.Range("D25").Value = "3.819.05"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -12.76%  "
$ws.Range("E28").Value = "  -14.55%  "
$ws.Range("D29").Value = "9.06"
$ws.Range("E29").Value = "  -12.29%  "
$ws.Range("E30").Value = "  -11.36%  "
$ws.Range("E31").Value = "  -15.59%  "
$ws.Range("D32").Value = "6.61"
$ws.Range("E32").Value = "  -10.33%  "
$ws.Range("E33").Value = "  -10.47%  "
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  -9.03%  "
$ws.Range("E36").Value = "  -8.32%  "
$ws.Range("E37").Value = "  -13.03%  "
$ws.Range("E38").Value = "  -7.84%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "2.22"
$ws.Range("E40").Value = "  -9.44%  "
$ws.Range("D41").Value = "0.0894"
$ws.Range("E41").Value = "  -11.72%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "0.937"
$ws.Range("E43").Value = "  -7.32%  "
$ws.Range("D44").Value = "164.97"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "47.71"
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").Value = "2.67"
$ws.Range("E46").Value = "  -17.73%  "
$ws.Range("E47").Value = "  -8.18%  "
$ws.Range("E50").Value = "  -7.07%  "
$ws.Range("D51").Value = "7.82"
$ws.Range("E51").Value = "  -9.31%  "

# Rows 48 and 49 swapped (InjectiveProtocol <-> FLOKI) with new values
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "0.000268"
$ws.Range("E48").Value = "  -12.69%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "27.59"
$ws.Range("E49").Value = "  -5.52%  "
